$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: round B5:AH5 to 2 decimal places (custom accuracy)
$newRow5 = @(15.37, 11.21, 1.07, 33.13, 27.41, 12.1, 45.18, 18.62, 8.2, 12.3, 13.38, 13.98, 3.86, 12.03, 17.06, 10.19, 0.83, 0.68, 175.79, 33.59, 11.1, 22.49, 12.02, 1.48, 22.03, 9.81, 8.75, 10.27, 13.99, 0.56, 40.69, 6.23, 13.88)

for ($i = 0; $i -lt $newRow5.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $newRow5[$i]
}

# Delete row 6 entirely (data trimmed)
$ws.Rows.Item(6).Delete()
